# Mise à jour des résultats du script
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear C207 (was "NA") back to an empty text cell ---
$c207 = $ws.Range("C207")
$c207.Value = "'"
$c207.Style = "Normal"

# --- Append new rows 208-221 ---
$newRows = @(
    @("2025-08-18", "réduction de la dérive", 1, 1),
    @("2025-08-18", "réduction de la dérive", 13, 2),
    @("2025-08-18", "zone tampon", 13, 1),
    @("2025-08-18", "buse", 14, 4),
    @("2025-08-18", "buse", 15, 2),
    @("2025-08-18", "buse", 16, 2),
    @("2025-08-18", "buse", 17, 2),
    @("2025-08-18", "buse", 18, 6),
    @("2025-08-18", "buse", 19, 1),
    @("2025-08-18", "buse", 20, 1),
    @("2025-08-18", "buse", 25, 3),
    @("2025-08-18", "zone tampon", 25, 5),
    @("2025-08-18", "herbicides", 26, 1),
    @("2025-08-18", "bonnes pratiques", 170, 2)
)

$startRow = 208
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]

    # Column A holds the date as plain text (not an Excel date value).
    # A leading apostrophe forces Excel to treat it as text instead of
    # auto-converting it to a date serial number; resetting the style
    # back to Normal afterwards drops the quote-prefix style so the
    # cell ends up with no style index, same as the source data.
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.Value = "'" + $data[0]
    $cellA.Style = "Normal"

    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
}
